# Updates cryptos.xlsx "Price" (D) and "Volume(1h)" (E) columns with
# refreshed values from the latest GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "64.109.70"
$ws.Range("E2").Value = "  +0.82%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.248.78"
$ws.Range("E3").Value = "  +3.95%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.63%  "

# Row 5: BNB
$ws.Range("D5").Value = "'586.55"
$ws.Range("E5").Value = "  -0.39%  "

# Row 6: Solana
$ws.Range("D6").Value = "'146.77"
$ws.Range("E6").Value = "  +0.41%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.66%  "

# Row 8: LidoStakedEther
$ws.Range("D8").Value = "3.148.28"
$ws.Range("E8").Value = "  +0.87%  "

# Row 9: XRP
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -1.17%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -0.38%  "

# Row 11: Toncoin
$ws.Range("D11").Value = "'5.83"
$ws.Range("E11").Value = "  +2.85%  "

# Row 12: Cardano
$ws.Range("E12").Value = "  -1.98%  "

# Row 13: ShibaInu
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  -2.35%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'37.06"
$ws.Range("E14").Value = "  +3.86%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.781.88"
$ws.Range("E15").Value = "  +3.78%  "

# Row 16: TRON
$ws.Range("E16").Value = "  -1.30%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "3.191.12"
$ws.Range("E17").Value = "  +1.95%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "63.965.32"
$ws.Range("E18").Value = "  +0.67%  "

# Row 19: Polkadot
$ws.Range("D19").Value = "'7.07"
$ws.Range("E19").Value = "  -1.05%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'464.81"
$ws.Range("E20").Value = "  +0.17%  "

# Row 21: Chainlink
$ws.Range("D21").Value = "'14.28"
$ws.Range("E21").Value = "  +0.84%  "

# Row 22: Polygon
$ws.Range("E22").Value = "  +0.07%  "

# Row 23: Uniswap
$ws.Range("D23").Value = "'7.43"
$ws.Range("E23").Value = "  -0.87%  "

# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").Value = "'12.93"
$ws.Range("E24").Value = "  -2.70%  "

# Row 25: Litecoin
$ws.Range("D25").Value = "'81.07"
$ws.Range("E25").Value = "  -0.80%  "

# Row 26: Fetch.AI
$ws.Range("E26").Value = "  +4.12%  "

# Row 27: Dai
$ws.Range("E27").Value = "  +0.27%  "

# Row 28: RenderToken
$ws.Range("D28").Value = "'9.27"
$ws.Range("E28").Value = "  +8.08%  "

# Row 29: FirstDigitalUSD
$ws.Range("E29").Value = "  +0.69%  "

# Row 30: PancakeSwap
$ws.Range("D30").Value = "'2.68"
$ws.Range("E30").Value = "  -0.66%  "

# Row 31: ImmutableX
$ws.Range("D31").Value = "'2.22"
$ws.Range("E31").Value = "  +0.97%  "

# Row 32: NEARProtocol
$ws.Range("D32").Value = "'7.10"
$ws.Range("E32").Value = "  +4.23%  "

# Row 33: EthereumClassic
$ws.Range("D33").Value = "'27.06"
$ws.Range("E33").Value = "  +0.68%  "

# Row 34: Hedera
$ws.Range("E34").Value = "  +0.69%  "

# Row 35: PEPE
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").Value = "  +0.23%  "

# Row 36: Mantle
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  -0.09%  "

# Row 37: Stacks
$ws.Range("E37").Value = "  -3.29%  "

# Row 38: dogwifhat
$ws.Range("D38").Value = "'3.31"
$ws.Range("E38").Value = "  -1.21%  "

# Row 39: Filecoin
$ws.Range("D39").Value = "'6.01"
$ws.Range("E39").Value = "  -1.75%  "

# Row 40: OKB
$ws.Range("D40").Value = "'51.43"
$ws.Range("E40").Value = "  +1.42%  "

# Row 41: Bittensor
$ws.Range("D41").Value = "'437.30"
$ws.Range("E41").Value = "  -2.65%  "

# Row 42: Cosmos
$ws.Range("E42").Value = "  +1.55%  "

# Row 43: VeChain
$ws.Range("E43").Value = "  -0.39%  "

# Row 44: Maker
$ws.Range("D44").Value = "2.911.65"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45: TheGraph
$ws.Range("E45").Value = "  +1.72%  "

# Row 46: Arweave
$ws.Range("D46").Value = "'39.05"
$ws.Range("E46").Value = "  +15.38%  "

# Row 47: Kaspa
$ws.Range("E47").Value = "  -3.47%  "

# Row 48: Monero
$ws.Range("D48").Value = "'126.51"
$ws.Range("E48").Value = "  -0.27%  "

# Row 50: Stellar
$ws.Range("E50").Value = "  -0.98%  "

# Row 51: ThetaToken
$ws.Range("E51").Value = "  +0.88%  "
